$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('UI_UX designer')
$ws.Range('C2').Value = 'Figma Essentials, Basic UI Elements, Components, Autolayout, and Variants, Typography Basics, Grids, Prototyping, UI Design, Clickable Prototyping'
$ws.Range('C3').Value = 'How to learn effectively, Figma practice'
$ws.Range('C4').Value = 'Intro, Design Process, Quantitative Research (Surveys), Heuristic Interface Analysis, Qualitative Research (In-Depth Interviews), Jobs To Be Done, Customer Journey Mapping, Kano Model. Features Prioritization, Information Architecture Mapping, User Flows, Prototyping, Unmoderated User Testing'
$ws.Range('C5').Value = 'Law of Proximity, Law of Common Region, Negative space, Contrast, Buttons, Inputs'
$ws.Range('C6').Value = 'Grids, Mobile Interfaces, Interface Patterns & Trends, Colors & UI, Typography, Design Gestalts and Rules, Handoff Preparation (UI Kits)'
$ws.Range('C7').Value = 'Scrum Methodology, Gathering Data From Interviews, Personas and User Stories, Product Hypotheses, Prototyping, Composition and Typography, Graphic Rhymes, Adaptive Design, How to Present Your Work, Project Presentation, Design Systems'
$ws.Range('C8').Value = 'Behance Overview, Case Structure, Trends, Case Design Tips, Second Behance Case'

$ws = $wb.Worksheets.Item('QA engineer')
$ws.Range('C2').Value = 'Introduction to SQL, SELECT Statement, WHERE Statement, NULL Value, LIKE, BETWEEN, and IN Statements, ORDER BY, LIMIT, DISTINCT, Aliases, Aggregate Functions, GROUP BY Statement, JOIN Statement, Functions, HAVING Statement'
$ws.Range('C3').Value = 'UI elements guide, Course Fundamentals'
$ws.Range('C4').Value = 'What is QA, Testing Types, What is a Bug Report, Bug Reports in details, Decomposition, Test Cases, Test Design Techniques, Test Checklist, Software Development Life Cycle, Software Testing Life Cycle, Test Plan, Requirements, TestRail & Jira'
$ws.Range('C5').Value = 'Environment Setup, Command Line Basics, Git Basics, Working With Branches, Working With Remote Repo (GitHub)'
$ws.Range('C6').Value = 'Git and Terminal'
$ws.Range('C7').Value = 'Introduction, Your First JavaScript Program, Main Concepts, Numbers, Strings, Boolean, Functions, Conditional Operators, Arrays, Loops, String Iteration, Strings Methods, Working With Arrays, Get Ready for the Interview'
$ws.Range('C8').Value = 'HTML Basics, CSS Basics, Colors and Fonts, Box Model Basics, Semantic Basics, Responsiveness Basics, CSS Selectors, Pseudo-Elements and Pseudo-Classes, Specificity, Links and URLs, Images, Media Queries, Forms, Position, Extra Topics'

$ws = $wb.Worksheets.Item('Python developer')
$ws.Range('C2').Value = 'Introduction, Main Concepts, Numbers, Strings, Boolean, Lists, Conditional Operators, Loops, Functions, Summary'
$ws.Range('C3').Value = 'How to learn effectively'
$ws.Range('C4').Value = 'Module Overview, Environment Setup, Code Style, Debugging, Working With Numbers, Working With Strings, Lists In Details, Dict Basics, Type Conversion, Loops In Details, Functions Revisited, Summary'
$ws.Range('C5').Value = 'Environment Setup, Code Style, Debugging, Working With Numbers, Working With Strings, Functions Revisited, Loops in Details, Type Conversion, Logical Operators, Lists in Details, Dict Basics, Extended Extra'
$ws.Range('C6').Value = 'Environment Setup, Command Line Basics, Git Basics, Working With Branches, Working With Remote Repo (GitHub)'
$ws.Range('C7').Value = 'How to Solve GitHub Tasks, Mutable Immutable Types, List and Dict Comprehensions, Functions in Details, Decorators, Classes, Classes in Details, Iterators and Generators, Modules and Imports, OOP Single Inheritance, OOP Multiple Inheritance, OOP Encapsulation Polymorphism Abstraction, Properties and Descriptors, Exception Handling, Exceptions in Details, File Handling, Memory Management, Testing, Testing in Details, Basic Modules Overview, Dict Advanced, Extra, Python Practice'
$ws.Range('C8').Value = 'Introduction to SQL, SELECT Statement, WHERE Statement, NULL Value, LIKE, BETWEEN, and IN Statements, ORDER BY, LIMIT, DISTINCT, Aliases, Aggregate Functions, GROUP BY Statement, JOIN Statement, Functions, HAVING Statement'

$ws = $wb.Worksheets.Item('Data analyst')
$ws.Range('C2').Value = 'Summary, Creating Dashboard, Calculated Fields, Netflix Cinematic Map, Introduction'
$ws.Range('C3').Value = 'Spreadsheet — the Basic Tool for Analytics'
$ws.Range('C4').Value = 'Introduction to SQL, SELECT Statement, WHERE Statement, NULL Value, LIKE, BETWEEN, and IN Statements, ORDER BY, LIMIT, DISTINCT, Aliases, Aggregate Functions, GROUP BY Statement, JOIN Statement, Functions, HAVING Statement'
$ws.Range('C5').Value = $null
$ws.Range('C6').Value = 'Introduction, Main Concepts, Numbers, Strings, Boolean, Lists, Conditional Operators, Loops, Functions, Summary'
$ws.Range('C7').Value = $null
$ws.Range('C8').Value = $null

$ws = $wb.Worksheets.Item('Digital marketer')
$ws.Range('C2').Value = 'Marketing Introduction, Offline Marketing vs. Online Marketing, Google Ads Search, Keywords, Your First Media Plan, Growth Opportunities, Course Overview'
$ws.Range('C3').Value = 'Market Sizing, Competitors Analysis, Customer Portrait, Points of Differentiation, Value Proposition'
$ws.Range('C4').Value = 'Google Ads Set Up, Campaign Structure, Campaign Objectives, Conversions, Campaign Types, Creating an Ad, Additional Settings of Google Ads Campaigns, Google Search Optimization'
$ws.Range('C5').Value = 'Campaign Creation, Targeting, Ads Creation, Google Display Ads Optimization'
$ws.Range('C6').Value = 'Meta Ads Overview, Meta Ads Setup, Meta Audience Types, Campaign Objectives, How to create Meta Ad Campaign, Creatives. Ad formats, Optimization of Advertising Campaigns'
$ws.Range('C7').Value = 'Introduction to SEO, Working with Keywords, On-Page Optimization, Working with Content, Link Building, What''s Next?, Competitor Analysis'
$ws.Range('C8').Value = 'Brand Awareness Channels, Outbound Lead Generation, Customer Activation, Portfolio Project, Customer Loyalty and Referral Programs'

$ws = $wb.Worksheets.Item('Front-end developer')
$ws.Range('C2').Value = 'Introduction, Your First JavaScript Program, Main Concepts, Numbers, Strings, Boolean, Functions, Conditional Operators, Arrays, Loops, String Iteration, Strings Methods, Working With Arrays, Get Ready for the Interview'
$ws.Range('C3').Value = 'Environment Setup, How to Learn Effectively'
$ws.Range('C4').Value = 'HTML Basics, CSS Basics, Colors and Fonts, Box Model Basics, Semantic Basics, Responsiveness Basics, CSS Selectors, Pseudo-Elements and Pseudo-Classes, Specificity, Links and URLs, Images, Media Queries, Forms, Position, Extra Topics'
$ws.Range('C5').Value = 'Environment Setup, Command Line Basics, Git Basics, Working With Branches, Working With Remote Repo (GitHub)'
$ws.Range('C6').Value = 'Code Style, Working With Numbers, Loops in Details, Working With Strings, Functions Revisited, Switch, Type Conversion, Logical Operators, Object Basics, Extended Extra'
$ws.Range('C7').Value = 'Environment Setup, Flexbox, BEM, Sass, Transformations and Animations, Grid, Landing [Portfolio Project], Document and Events'
$ws.Range('C8').Value = 'How to Solve Tasks on Github, Object Advanced, Methods, Array Methods, Callbacks, Array Iteration Methods Implementation, Array Iteration Methods Usage, Array Iteration Methods Practice, Closures, JS Practice, Prototype, Constructors, Classes, Extra'

$ws = $wb.Worksheets.Item('Full-stack developer')
$ws.Range('C2').Value = 'Introduction, Your First JavaScript Program, Main Concepts, Numbers, Strings, Boolean, Functions, Conditional Operators, Arrays, Loops, String Iteration, Strings Methods, Working With Arrays, Get Ready for the Interview'
$ws.Range('C3').Value = 'How to Learn Effectively, Environment Setup'
$ws.Range('C4').Value = 'HTML Basics, CSS Basics, Colors and Fonts, Box Model Basics, Semantic Basics, Responsiveness Basics, CSS Selectors, Pseudo-Elements and Pseudo-Classes, Specificity, Links and URLs, Images, Media Queries, Forms, Position, Extra Topics'
$ws.Range('C5').Value = 'Environment Setup, Command Line Basics, Git Basics, Working With Branches, Working With Remote Repo (GitHub)'
$ws.Range('C6').Value = 'Code Style, Working With Numbers, Loops in Details, Working With Strings, Functions Revisited, Switch, Type Conversion, Logical Operators, Object Basics, Extended Extra'
$ws.Range('C7').Value = 'Environment Setup, Flexbox, BEM, Sass, Transformations and Animations, Grid, Landing [Portfolio Project], Document and Events'
$ws.Range('C8').Value = 'How to Solve Tasks on Github, Object Advanced, Methods, Array Methods, Callbacks, Array Iteration Methods Implementation, Array Iteration Methods Usage, Array Iteration Methods Practice, Closures, JS Practice, Prototype, Constructors, Classes, Extra'

$ws = $wb.Worksheets.Item('DevOps engineer')
$ws.Range('C2').Value = 'Introduction, Main Concepts, Numbers, Strings, Boolean, Lists, Conditional Operators, Loops, Functions, Summary'
$ws.Range('C3').Value = 'Module Overview, Environment Setup, Code Style, Debugging, Working With Numbers, Working With Strings, Lists In Details, Dict Basics, Type Conversion, Loops In Details, Functions Revisited, Summary'
$ws.Range('C4').Value = 'Environment Setup, Code Style, Debugging, Working With Numbers, Working With Strings, Functions Revisited, Loops in Details, Type Conversion, Logical Operators, Lists in Details, Dict Basics, Extended Extra'
$ws.Range('C5').Value = 'Environment Setup, Command Line Basics, Git Basics, Working With Branches, Working With Remote Repo (GitHub)'
$ws.Range('C6').Value = 'Introduction, DevOps Culture and Practices, DevOps as a Role, SDLC and Development Methodologies, Summary'
$ws.Range('C7').Value = 'Introduction to the Web, OSI Model, HTTP, API Interface, Encryption, Security in the Web'
$ws.Range('C8').Value = 'Environment Setup, What Is Operating System, Working With Files, Processes and Services, Resource Management, Managing Users, Installing Software, Using SSH, Shell Scripting, Web Server'

$ws = $wb.Worksheets.Item('Java developer')
$ws.Range('C2').Value = 'Introduction, Main Concepts, Data Types, Operators, Methods, Boolean, Conditionals, Loops, Arrays, String, Practice'
$ws.Range('C3').Value = 'Type Casting, Arrays Extended, Switch Case, Classes, Classes Advanced, Methods Advanced, Practice'
$ws.Range('C4').Value = 'Environment Setup, Command Line Basics, Git Basics, Working With Branches, Working With Remote Repo (GitHub)'
$ws.Range('C5').Value = 'Intro, Java and IntelliJ Idea, Maven, Checkstyle and Review process'
$ws.Range('C6').Value = 'How to Solve Tasks, Memory, Wrappers, String, Bit Manipulation, OOP, Abstract Class vs. Interface, Exception, Files, Equals, Hashcode, and Clone, Patterns and Recursion, Generics, ArrayList, LinkedList, PECS, HashMap, Array Practice, Immutable, Set, Queue, Stack, and Comparator, Java JUnit, Java 8 — Part 1, Java 8 — Part 2, Java 8 — Part 3, Java 9-19 Improvements, Stream API Practice, Java SOLID, Java Dependency Injection, Logger'
$ws.Range('C7').Value = 'Introduction to SQL, SELECT Statement, WHERE Statement, NULL Value, LIKE, BETWEEN, and IN Statements, ORDER BY, LIMIT, DISTINCT, Aliases, Aggregate Functions, GROUP BY Statement, JOIN Statement, Functions, HAVING Statement'
$ws.Range('C8').Value = 'Database Intro, JDBC Intro, Hibernate Entities and Transactions, Entity Relations and Life Cycle, JPA and Cascade Types, Working With Lazy Initialization, Criteria API, Transaction Isolation Levels, N+1 Problem, Database Change Management Tools, SQL Practice, Hibernate Inheritance, Hibernate Cache Levels, Performance Improvement'

$ws = $wb.Worksheets.Item('Recruiter')
$ws.Range('C2').Value = 'General Info About IT, Company Types and UA Market Overview, Non-Tech Roles Overview, Tech Roles Overview'
$ws.Range('C3').Value = 'Software Development Life Cycle, Technologies'
$ws.Range('C4').Value = 'Recruitment Process Overview, Roles in the Recruitment Process, Vacancy Creation and Posting, Profile Screening'
$ws.Range('C5').Value = 'Boolean Search, X-Ray, LinkedIn, GitHub, Stack Overflow, Kaggle, Dou, Meetup, Other Resources for Sourcing, Mails & Communication, Sourcing Strategy & Candidate Database Management, Sourcing Strategy for CV'
$ws.Range('C6').Value = 'Introduction, Preparing Your Resume, Interview Preparation, Preparing Your Cover Letter, DOU and Google Profiles, Creating Your LinkedIn Profile, How to Communicate in the Right Way, Creating Your Djinni Profile, How to Find Vacancies, How to Apply for a Job on Different Platforms, How to Find Recruiters'' Contacts, Employment Daily Activities, Work Details FAQ'
$ws.Range('C7').Value = 'Interview & Prescreening Structure, Feedback & Job Offer, Recruitment Funnel'
$ws.Range('C8').Value = 'Introduction to the Web, OSI Model, HTTP, API Interface, Encryption, Security in the Web'
